$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all data rows (2-43)
# from serial date 45828 (2025-06-20) to 45829 (2025-06-21)
for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45828) {
        $cell.Value = 45829
    }
}
